$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy header formatting from H1 into the two new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$iValues = @(7,7,7,8,8,8,9,9,8,6,6,6,7,8,8,8,6,9,8,6,6,7,10,9,6,6,5,7,6,5,5,6,6,9,8,8,8,6,8,8,9,9,7,9,9,9,8,9,8,8,9,8,9,7,10,8,9,7,9,9,7,8,9,6,8,7,6,7,7,6,6,6,8,8,7,10,7,8,5,5,7,7,7,6,5,7,4,3,3)
$jValues = @(7,7,7,8,8,8,9,9,8,6,6,6,7,8,8,8,6,9,8,6,6,7,10,9,6,6,6,7,7,5,5,6,6,9,8,8,8,7,8,9,9,9,7,9,9,9,8,9,8,8,9,9,9,7,10,8,9,7,9,9,7,8,9,6,8,7,7,7,7,7,6,6,8,8,7,10,7,8,5,5,7,7,7,6,5,7,4,3,3)

$rowCount = $iValues.Length
$data = New-Object 'object[,]' $rowCount,2
for ($k = 0; $k -lt $rowCount; $k++) {
    $data[$k,0] = $iValues[$k]
    $data[$k,1] = $jValues[$k]
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$rangeAddr = "I" + $startRow + ":J" + $endRow
$ws.Range($rangeAddr).Value2 = $data
